$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.222.84'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '2.423.31'
$ws.Range("E3").Value = '  -2.32%  '
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = $origStyle
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.92'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +0.32%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.70'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -4.41%  '
$ws.Range("E7").Value = '  -0.01%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.525'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -2.83%  '
$ws.Range("D9").Value = '2.418.51'
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("E10").Value = '  -5.30%  '
$ws.Range("E11").Value = '  +0.86%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.18'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  -1.80%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.345'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -3.87%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.30'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -3.36%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000172'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -4.93%  '
$ws.Range("D16").Value = '2.869.85'
$ws.Range("E16").Value = '  -2.37%  '
$ws.Range("D17").Value = '62.091.14'
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").Value = '2.428.47'
$ws.Range("E18").Value = '  -2.56%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.94'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -5.16%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.07'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  -3.99%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.38'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -0.17%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.11'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -2.17%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.98'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -5.34%  '
$ws.Range("E24").Value = '  +0.17%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.56'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -0.86%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '629.39'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  +0.65%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.96'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +4.02%  '
$ws.Range("D28").Value = '2.544.08'
$ws.Range("E28").Value = '  -2.88%  '
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D30").Value = '0.0₃0943'
$ws.Range("E30").Value = '  -9.46%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -7.07%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.00'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -4.59%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.88'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("E34").Value = '  -3.43%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.94'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -5.72%  '
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("E37").Value = '  -6.08%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.373'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '148.09'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.32'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -2.65%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.23'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -4.36%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.74'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -4.75%  '
$ws.Range("E43").Value = '  -0.01%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.19'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +0.70%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -9.32%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '143.07'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -4.41%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.65'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -2.96%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0518'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -4.95%  '
$ws.Range("E49").Value = '  -2.40%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.40'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -9.09%  '
$ws.Range("E51").Value = '  +7.37%  '
